# New crime data collected - update the 84th Precinct weekly CompStat sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: volume/number line and the "report covering the week" line.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/27/2023  Through  4/2/2023"

# ---------------------------------------------------------------------------
# Helper: copy number-format/style from a "template" cell onto a target cell,
# then drop in the new value (keeps General vs #,##0 vs text styling intact
# exactly like Excel's Home > Format Painter would).
# ---------------------------------------------------------------------------
function Set-StyledValue($targetAddr, $templateAddr, $value) {
    $ws.Range($templateAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($targetAddr).Value = $value
}

function Set-StyledText($targetAddr, $templateAddr, $value) {
    $ws.Range($templateAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($templateAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4163) | Out-Null   # xlPasteValues (keeps as shared-string text)
}

# ---------------------------------------------------------------------------
# Row 14 — Murder
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -75

# ---------------------------------------------------------------------------
# Row 15 — Rape (28-day columns flip from numbers to the "N/A" placeholders)
# ---------------------------------------------------------------------------
Set-StyledText "G15" "C14" "0"
Set-StyledText "H15" "E14" "***.*"

# ---------------------------------------------------------------------------
# Row 16 — Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 24
$ws.Range("J16").Value = 35
$ws.Range("K16").Value = -31.428571428571
$ws.Range("L16").Value = -25
$ws.Range("M16").Value = -58.620689655172
$ws.Range("N16").Value = -91.808873720136

# ---------------------------------------------------------------------------
# Row 17 — Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 175
$ws.Range("I17").Value = 66
$ws.Range("J17").Value = 37
$ws.Range("K17").Value = 78.378378378378
$ws.Range("L17").Value = 112.903225806452
$ws.Range("M17").Value = 106.25
$ws.Range("N17").Value = -34.653465346534

# ---------------------------------------------------------------------------
# Row 18 — Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 59
$ws.Range("K18").Value = -3.389830508474
$ws.Range("L18").Value = 46.153846153846
$ws.Range("M18").Value = 185
$ws.Range("N18").Value = -67.428571428571

# ---------------------------------------------------------------------------
# Row 19 — Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 140
$ws.Range("F19").Value = 60
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 20
$ws.Range("I19").Value = 167
$ws.Range("J19").Value = 153
$ws.Range("K19").Value = 9.150326797385
$ws.Range("L19").Value = 38.016528925619
$ws.Range("M19").Value = 59.047619047619
$ws.Range("N19").Value = -33.992094861660

# ---------------------------------------------------------------------------
# Row 20 — G.L.A. (week-to-date columns flip from "N/A" placeholders to numbers)
# ---------------------------------------------------------------------------
Set-StyledValue "C20" "C16" 2
Set-StyledValue "D20" "D16" 1
Set-StyledValue "E20" "E16" 100
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 14
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = 7.692307692307
$ws.Range("L20").Value = 133.333333333333
$ws.Range("M20").Value = -6.666666666666
$ws.Range("N20").Value = -92.265193370165

# ---------------------------------------------------------------------------
# Row 21 — TOTAL (bold summary row)
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = 68
$ws.Range("F21").Value = 112
$ws.Range("H21").Value = 20.430107526881
$ws.Range("I21").Value = 331
$ws.Range("J21").Value = 301
$ws.Range("K21").Value = 9.966777408637
$ws.Range("L21").Value = 42.672413793103
$ws.Range("M21").Value = 43.290043290043
$ws.Range("N21").Value = -67.195242814668

# ---------------------------------------------------------------------------
# Row 22 — Transit
# ---------------------------------------------------------------------------
Set-StyledValue "D22" "C22" 1
Set-StyledValue "E22" "H22" 200
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 133.333333333333
$ws.Range("I22").Value = 12
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = 9.090909090909
$ws.Range("L22").Value = 9.090909090909
$ws.Range("M22").Value = -29.411764705882

# ---------------------------------------------------------------------------
# Row 23 — Housing
# ---------------------------------------------------------------------------
Set-StyledValue "C23" "C22" 2
Set-StyledValue "D23" "C22" 1
Set-StyledValue "E23" "H22" 100
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -40
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 12
$ws.Range("K23").Value = -8.333333333333
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 175

# ---------------------------------------------------------------------------
# Row 24 — Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -5.128205128205
$ws.Range("G24").Value = 158
$ws.Range("H24").Value = -4.430379746835
$ws.Range("I24").Value = 478
$ws.Range("J24").Value = 378
$ws.Range("K24").Value = 26.455026455026
$ws.Range("L24").Value = 43.113772455089
$ws.Range("M24").Value = 43.543543543543

# ---------------------------------------------------------------------------
# Row 25 — Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -70
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 25
$ws.Range("H25").Value = 8
$ws.Range("I25").Value = 93
$ws.Range("J25").Value = 81
$ws.Range("K25").Value = 14.814814814814
$ws.Range("L25").Value = 82.352941176470
$ws.Range("M25").Value = -11.428571428571

# ---------------------------------------------------------------------------
# Row 26 — UCR Rape* (28-day columns flip from numbers to "N/A" placeholders)
# ---------------------------------------------------------------------------
Set-StyledText "G26" "C14" "0"
Set-StyledText "H26" "E14" "***.*"

# ---------------------------------------------------------------------------
# Row 27 — Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -25
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -30
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -22.222222222222
$ws.Range("L27").Value = 7.692307692307

# ---------------------------------------------------------------------------
# Row 28 — Shooting Vic.
# ---------------------------------------------------------------------------
$ws.Range("N28").Value = -72.727272727272

# ---------------------------------------------------------------------------
# Row 29 — Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Range("N29").Value = -70

# ---------------------------------------------------------------------------
# Row 30 — Hate Crimes
# ---------------------------------------------------------------------------
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 100
